$wb = $excel.ActiveWorkbook

# --- Text change: "Ready for handoff" -> "In Translation" -------------------
# This status string appears in the Overview sheet (columns "zh-cn"/"de-de",
# row 2) and in each per-locale sheet's "Status" column (C2).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- Column width changes ----------------------------------------------------
# Overview!E:F and zh-cn!C / de-de!C shrink from ~17.216 to ~13.410 characters.
# Excel quantizes ColumnWidth to whole-pixel increments (1/6 of a character
# unit here) when it is written back out, and stores the width as
# round(6*ColumnWidth)/6 + 5/6. To land as close as possible on the target
# stored width we therefore feed it "target - 5/6" rather than the target
# value itself.
$targetStoredWidth = 13.4101845877511
$columnWidthInput = $targetStoredWidth - (5 / 6)

$wsOverview.Columns.Item(5).ColumnWidth = $columnWidthInput   # column E
$wsOverview.Columns.Item(6).ColumnWidth = $columnWidthInput   # column F

$wsZhCn.Columns.Item(3).ColumnWidth = $columnWidthInput       # column C
$wsDeDe.Columns.Item(3).ColumnWidth = $columnWidthInput       # column C
